$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$reqFermentativos = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$reqBebidas = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)`n"

# Swap order: LOT2052 (Bebidas) now appears first (row 23), LOT2028 (Fermentativos) second (row 24)
$ws.Range("B23").Value = $reqBebidas
$ws.Range("C23").Value = $reqBebidas

$ws.Range("B24").Value = $reqFermentativos
$ws.Range("C24").Value = $reqFermentativos
